$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2 (CasesTab query): append an ORDER BY / LIMIT clause -----------------
$b2 = $ws.Range("B2").Formula
$b2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"
$ws.Range("B2").Formula = $b2

# --- B3 (SamplesTab query): append an ORDER BY / LIMIT clause --------------
$b3 = $ws.Range("B3").Formula
$b3 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Formula = $b3

# --- B4 (FilesTab query): swap the trailing "order by" for the new clause --
$b4 = $ws.Range("B4").Formula
$b4 = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Formula = $b4

# --- Row heights grew because the wrapped text got longer ------------------
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# --- Selection / scroll position (best effort UI state) --------------------
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 3 } catch {}
$ws.Range("B4").Select()
